$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
# Row 38 (ALC)
$ws_ALC.Cells.Item(38, 8).Value = 865.8333
$ws_ALC.Cells.Item(38, 9).Value = 439
$ws_ALC.Cells.Item(38, 10).Value = 3000
$ws_ALC.Cells.Item(38, 11).Value = 1317
$ws_ALC.Cells.Item(38, 12).Value = 9000
$ws_ALC.Cells.Item(38, 13).Value = -945
$ws_ALC.Cells.Item(38, 14).Value = -9744

# Row 69 (ALC)
$ws_ALC.Cells.Item(69, 8).Value = 3799.75
$ws_ALC.Cells.Item(69, 9).Value = 2733
$ws_ALC.Cells.Item(69, 11).Value = 8199
$ws_ALC.Cells.Item(69, 13).Value = -7325

# Row 72 (ALC)
$ws_ALC.Cells.Item(72, 8).Value = 3799.75
$ws_ALC.Cells.Item(72, 9).Value = 2733
$ws_ALC.Cells.Item(72, 11).Value = 24597
$ws_ALC.Cells.Item(72, 13).Value = -20229

# Row 131 (ALC)
$ws_ALC.Cells.Item(131, 8).Value = 1994.4
$ws_ALC.Cells.Item(131, 9).Value = 977
$ws_ALC.Cells.Item(131, 11).Value = 2931
$ws_ALC.Cells.Item(131, 13).Value = 2109

$ws_ARM = $wb.Worksheets.Item("ARM")
# Row 31 (ARM)
$ws_ARM.Cells.Item(31, 8).Value = 4704.5
$ws_ARM.Cells.Item(31, 9).Value = 4704.5
$ws_ARM.Cells.Item(31, 11).Value = 4704.5
$ws_ARM.Cells.Item(31, 13).Value = -4410.5

# Row 61 (ARM)
$ws_ARM.Cells.Item(61, 8).Value = 2598.9092
$ws_ARM.Cells.Item(61, 9).Value = 1357.875
$ws_ARM.Cells.Item(61, 10).Value = 5908.3335
$ws_ARM.Cells.Item(61, 11).Value = 1357.875
$ws_ARM.Cells.Item(61, 12).Value = 5908.3335
$ws_ARM.Cells.Item(61, 13).Value = -1145.875
$ws_ARM.Cells.Item(61, 14).Value = -6332.3335

# Row 122 (ARM)
$ws_ARM.Cells.Item(122, 8).Value = 2858.375
$ws_ARM.Cells.Item(122, 9).Value = 1838.1428
$ws_ARM.Cells.Item(122, 11).Value = 5514.428400000001
$ws_ARM.Cells.Item(122, 13).Value = -3064.428400000001

# Row 132 (ARM)
$ws_ARM.Cells.Item(132, 8).Value = 1793.4642
$ws_ARM.Cells.Item(132, 9).Value = 1196.0555
$ws_ARM.Cells.Item(132, 11).Value = 3588.1665
$ws_ARM.Cells.Item(132, 13).Value = -1058.1665

# Row 136 (ARM)
$ws_ARM.Cells.Item(136, 8).Value = 2598.9092
$ws_ARM.Cells.Item(136, 9).Value = 1357.875
$ws_ARM.Cells.Item(136, 10).Value = 5908.3335
$ws_ARM.Cells.Item(136, 11).Value = 4073.625
$ws_ARM.Cells.Item(136, 12).Value = 17725.0005
$ws_ARM.Cells.Item(136, 13).Value = -1523.625
$ws_ARM.Cells.Item(136, 14).Value = -22825.0005

$ws_BSM = $wb.Worksheets.Item("BSM")
# Row 99 (BSM)
$ws_BSM.Cells.Item(99, 8).Value = 1450
$ws_BSM.Cells.Item(99, 9).Value = 1209.091
$ws_BSM.Cells.Item(99, 11).Value = 1209.091
$ws_BSM.Cells.Item(99, 13).Value = 288.9090000000001

# Row 122 (BSM)
$ws_BSM.Cells.Item(122, 8).Value = 68000
$ws_BSM.Cells.Item(122, 10).Value = 68000
$ws_BSM.Cells.Item(122, 12).Value = 68000
$ws_BSM.Cells.Item(122, 14).Value = -77800

# Row 134 (BSM)
$ws_BSM.Cells.Item(134, 8).Value = 3554.4807
$ws_BSM.Cells.Item(134, 9).Value = 3902.0889
$ws_BSM.Cells.Item(134, 10).Value = 1319.8572
$ws_BSM.Cells.Item(134, 11).Value = 11706.2667
$ws_BSM.Cells.Item(134, 12).Value = 3959.5716
$ws_BSM.Cells.Item(134, 13).Value = -9171.2667
$ws_BSM.Cells.Item(134, 14).Value = -9029.571599999999

$ws_CRP = $wb.Worksheets.Item("CRP")
# Row 16 (CRP)
$ws_CRP.Cells.Item(16, 8).Value = 836.25
$ws_CRP.Cells.Item(16, 9).Value = 777.3333
$ws_CRP.Cells.Item(16, 11).Value = 777.3333
$ws_CRP.Cells.Item(16, 13).Value = -490.3333

# Row 31 (CRP)
$ws_CRP.Cells.Item(31, 8).Value = 1321.1428
$ws_CRP.Cells.Item(31, 9).Value = 1122.6818
$ws_CRP.Cells.Item(31, 10).Value = 2048.8333
$ws_CRP.Cells.Item(31, 11).Value = 1122.6818
$ws_CRP.Cells.Item(31, 12).Value = 2048.8333
$ws_CRP.Cells.Item(31, 13).Value = -827.6818000000001
$ws_CRP.Cells.Item(31, 14).Value = -2638.8333

# Row 34 (CRP)
$ws_CRP.Cells.Item(34, 8).Value = 1321.1428
$ws_CRP.Cells.Item(34, 9).Value = 1122.6818
$ws_CRP.Cells.Item(34, 10).Value = 2048.8333
$ws_CRP.Cells.Item(34, 11).Value = 1122.6818
$ws_CRP.Cells.Item(34, 12).Value = 2048.8333
$ws_CRP.Cells.Item(34, 13).Value = -920.6818000000001
$ws_CRP.Cells.Item(34, 14).Value = -2452.8333

# Row 113 (CRP)
$ws_CRP.Cells.Item(113, 8).Value = 836.25
$ws_CRP.Cells.Item(113, 9).Value = 777.3333
$ws_CRP.Cells.Item(113, 11).Value = 777.3333
$ws_CRP.Cells.Item(113, 13).Value = 1392.6667

# Row 132 (CRP)
$ws_CRP.Cells.Item(132, 8).Value = 1735.5272
$ws_CRP.Cells.Item(132, 9).Value = 1178.973
$ws_CRP.Cells.Item(132, 11).Value = 3536.919
$ws_CRP.Cells.Item(132, 13).Value = -1006.919

# Row 141 (CRP)
$ws_CRP.Cells.Item(141, 8).Value = 71539
$ws_CRP.Cells.Item(141, 10).Value = 71539
$ws_CRP.Cells.Item(141, 12).Value = 71539
$ws_CRP.Cells.Item(141, 14).Value = -81899

$ws_CUL = $wb.Worksheets.Item("CUL")
# Row 98 (CUL)
$ws_CUL.Cells.Item(98, 8).Value = 316.5
$ws_CUL.Cells.Item(98, 9).Value = 225
$ws_CUL.Cells.Item(98, 10).Value = 362.25
$ws_CUL.Cells.Item(98, 11).Value = 675
$ws_CUL.Cells.Item(98, 12).Value = 1086.75
$ws_CUL.Cells.Item(98, 13).Value = 823
$ws_CUL.Cells.Item(98, 14).Value = -4082.75

# Row 114 (CUL)
$ws_CUL.Cells.Item(114, 8).Value = 17858932
$ws_CUL.Cells.Item(114, 9).Value = 864
$ws_CUL.Cells.Item(114, 10).Value = 23811622
$ws_CUL.Cells.Item(114, 11).Value = 2592
$ws_CUL.Cells.Item(114, 12).Value = 71434866
$ws_CUL.Cells.Item(114, 13).Value = 662
$ws_CUL.Cells.Item(114, 14).Value = -71441374

# Row 131 (CUL)
$ws_CUL.Cells.Item(131, 8).Value = 801.5
$ws_CUL.Cells.Item(131, 9).Value = 410.77777
$ws_CUL.Cells.Item(131, 10).Value = 840.1429000000001
$ws_CUL.Cells.Item(131, 11).Value = 1232.33331
$ws_CUL.Cells.Item(131, 12).Value = 2520.4287
$ws_CUL.Cells.Item(131, 13).Value = 3807.66669
$ws_CUL.Cells.Item(131, 14).Value = -12600.4287

# Row 134 (CUL)
$ws_CUL.Cells.Item(134, 8).Value = 2352.2307
$ws_CUL.Cells.Item(134, 9).Value = 1858.8
$ws_CUL.Cells.Item(134, 11).Value = 5576.4
$ws_CUL.Cells.Item(134, 13).Value = -506.3999999999996

$ws_GSM = $wb.Worksheets.Item("GSM")
# Row 21 (GSM)
$ws_GSM.Cells.Item(21, 8).Value = 8380171.5
$ws_GSM.Cells.Item(21, 9).Value = 25000500
$ws_GSM.Cells.Item(21, 11).Value = 25000500
$ws_GSM.Cells.Item(21, 13).Value = -25000327

# Row 30 (GSM)
$ws_GSM.Cells.Item(30, 8).Value = 8380171.5
$ws_GSM.Cells.Item(30, 9).Value = 25000500
$ws_GSM.Cells.Item(30, 11).Value = 25000500
$ws_GSM.Cells.Item(30, 13).Value = -25000395

# Row 122 (GSM)
$ws_GSM.Cells.Item(122, 8).Value = 1546.9375
$ws_GSM.Cells.Item(122, 10).Value = 1929.8125
$ws_GSM.Cells.Item(122, 12).Value = 5789.4375
$ws_GSM.Cells.Item(122, 14).Value = -10689.4375

$ws_LTW = $wb.Worksheets.Item("LTW")
# Row 26 (LTW)
$ws_LTW.Cells.Item(26, 8).Value = 0
$ws_LTW.Cells.Item(26, 9).Value = 0
$ws_LTW.Cells.Item(26, 11).Value = 0
$ws_LTW.Cells.Item(26, 13).ClearContents()

# Row 40 (LTW)
$ws_LTW.Cells.Item(40, 8).Value = 3973.8235
$ws_LTW.Cells.Item(40, 9).Value = 1564.5834
$ws_LTW.Cells.Item(40, 11).Value = 1564.5834
$ws_LTW.Cells.Item(40, 13).Value = -1428.5834

# Row 122 (LTW)
$ws_LTW.Cells.Item(122, 8).Value = 5973
$ws_LTW.Cells.Item(122, 10).Value = 9174.75
$ws_LTW.Cells.Item(122, 12).Value = 27524.25
$ws_LTW.Cells.Item(122, 14).Value = -32424.25

$ws_WVR = $wb.Worksheets.Item("WVR")
# Row 96 (WVR)
$ws_WVR.Cells.Item(96, 8).Value = 4000
$ws_WVR.Cells.Item(96, 9).Value = 0
$ws_WVR.Cells.Item(96, 11).Value = 0
$ws_WVR.Cells.Item(96, 13).ClearContents()

# Row 113 (WVR)
$ws_WVR.Cells.Item(113, 8).Value = 1155.7778
$ws_WVR.Cells.Item(113, 9).Value = 960.2
$ws_WVR.Cells.Item(113, 10).Value = 1400.25
$ws_WVR.Cells.Item(113, 11).Value = 2880.6
$ws_WVR.Cells.Item(113, 12).Value = 4200.75
$ws_WVR.Cells.Item(113, 13).Value = -710.6000000000004
$ws_WVR.Cells.Item(113, 14).Value = -8540.75

# Row 122 (WVR)
$ws_WVR.Cells.Item(122, 8).Value = 61036
$ws_WVR.Cells.Item(122, 9).Value = 61036
$ws_WVR.Cells.Item(122, 11).Value = 183108
$ws_WVR.Cells.Item(122, 13).Value = -180658

# Row 132 (WVR)
$ws_WVR.Cells.Item(132, 8).Value = 1169.079
$ws_WVR.Cells.Item(132, 9).Value = 852.3226
$ws_WVR.Cells.Item(132, 10).Value = 2571.8572
$ws_WVR.Cells.Item(132, 11).Value = 2556.9678
$ws_WVR.Cells.Item(132, 12).Value = 7715.571599999999
$ws_WVR.Cells.Item(132, 13).Value = -26.9677999999999
$ws_WVR.Cells.Item(132, 14).Value = -12775.5716
